$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.323.96"
$ws.Range("E2").Value = "  -3.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.936.70"
$ws.Range("E3").Value = "  -3.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.40"
$ws.Range("E5").Value = "  -2.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7237"
$ws.Range("E6").Value = "  -6.71%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3313"
$ws.Range("E8").Value = "  -4.66%  "

$ws.Range("E9").Value = "  -1.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07288"
$ws.Range("E10").Value = "  +1.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8088"
$ws.Range("E11").Value = "  -4.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08086"
$ws.Range("E12").Value = "  -1.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.934.19"
$ws.Range("E13").Value = "  -3.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.492"
$ws.Range("E14").Value = "  -2.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.57"
$ws.Range("E15").Value = "  -6.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.11"
$ws.Range("E16").Value = "  -3.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.318.66"
$ws.Range("E17").Value = "  -3.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008343"
$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "250.63"
$ws.Range("E19").Value = "  -8.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.875"
$ws.Range("E20").Value = "  -2.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.188.99"
$ws.Range("E21").Value = "  -3.24%  "

$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.990"
$ws.Range("E24").Value = "  -2.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.758"
$ws.Range("E25").Value = "  -3.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.05"
$ws.Range("E26").Value = "  -0.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.381"
$ws.Range("E27").Value = "  -1.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.27"
$ws.Range("E28").Value = "  -3.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1327"
$ws.Range("E29").Value = "  -6.41%  "

$ws.Range("E30").Value = "  -2.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.347"
$ws.Range("E31").Value = "  -1.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.425"
$ws.Range("E32").Value = "  -4.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.177"
$ws.Range("E33").Value = "  -6.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05197"
$ws.Range("E34").Value = "  -3.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.285"
$ws.Range("E35").Value = "  +1.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7489"
$ws.Range("E36").Value = "  -5.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.747"
$ws.Range("E37").Value = "  -1.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01976"
$ws.Range("E38").Value = "  -1.89%  "

$ws.Range("E39").Value = "  -3.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "79.08"
$ws.Range("E40").Value = "  -7.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.349"
$ws.Range("E41").Value = "  -7.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4530"
$ws.Range("E42").Value = "  -3.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.023"
$ws.Range("E43").Value = "  -5.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8475"
$ws.Range("E44").Value = "  -1.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.0000"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.77"
$ws.Range("E46").Value = "  -3.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.719"
$ws.Range("E47").Value = "  -4.71%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.458"
$ws.Range("E48").Value = "  -4.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.71"
$ws.Range("E49").Value = "  -3.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4186"
$ws.Range("E50").Value = "  -4.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06035"
$ws.Range("E51").Value = "  -0.44%  "
